$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "aaa"
$ws.Range("B4").Value = "aaa"
$ws.Range("C4").Value = "Failed! please enter strong password"

$ws.Range("C4").Select()
